# Add a new "Count" column (I) to the Shop sheet, mirroring the existing
# per-row semantics used by the other metadata columns (B..H):
#   - Row 1   : header label "Count"
#   - Row 2   : the column's declared type ("int", same as the other data
#               columns use in this type row)
#   - Rows 3-6: Public/Private/Save/Cache boolean flags, copied from column H
#   - Rows 7-8: Ref/Upload boolean flags, copied from column H
#   - Row 9   : left untouched (section header row, no data columns)
#   - Rows 10-71: the actual per-item data, always 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "Count"

# Type marker row - same as the rest of the row (int), copy formatting from
# column H too so the new cell matches the existing B:H block exactly.
$srcI2 = $ws.Cells.Item(2, 8)
$dstI2 = $ws.Cells.Item(2, 9)
$dstI2.Value = $srcI2.Value()
$dstI2.HorizontalAlignment = $srcI2.HorizontalAlignment
$dstI2.WrapText = $srcI2.WrapText
$dstI2.Borders.LineStyle = $srcI2.Borders.LineStyle
$dstI2.Interior.Color = $srcI2.Interior.Color

# Boolean flag rows (Public, Private, Save, Cache) - copy value + formatting
# from column H so the new cells match the existing B:H block exactly.
foreach ($r in 3..6) {
    $srcCell = $ws.Cells.Item($r, 8)
    $dstCell = $ws.Cells.Item($r, 9)
    $dstCell.Value = $srcCell.Value()
    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $dstCell.WrapText = $srcCell.WrapText
    $dstCell.Borders.LineStyle = $srcCell.Borders.LineStyle
    $dstCell.Interior.Color = $srcCell.Interior.Color
}

# Ref / Upload boolean flag rows - same value as column H (formatting here
# already matches the row default, so a plain value write is enough).
foreach ($r in 7..8) {
    $srcCell = $ws.Cells.Item($r, 8)
    $dstCell = $ws.Cells.Item($r, 9)
    $dstCell.Value = $srcCell.Value()
}

# Data rows: every shop entry gets a Count of 1.
foreach ($r in 10..71) {
    $ws.Cells.Item($r, 9).Value = 1
}
